$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.000.21'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.012.75'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +3.05%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.68'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.06'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -1.60%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.390'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +3.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0800'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +1.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("E12").Value = '  +5.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.307.96'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +3.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.842'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.05'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +2.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.41'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +2.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.010.70'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +3.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.979.89'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.01'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.19'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.31'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +5.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.31'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.48'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.93%  '

$ws.Range("E28").Value = '  -4.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.68'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +2.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +6.81%  '

$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0667'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +9.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.55'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +13.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.45'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.60'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +6.29%  '

$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("E38").Value = '  +1.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.33'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("E40").Value = '  +3.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0966'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0215'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +3.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.17'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.48'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +4.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.88'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +2.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.370.09'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +0.65%  '

$ws.Range("E47").Value = '  +2.07%  '

$ws.Range("E48").Value = '  +3.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.13'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +15.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.87'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.91'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.71%  '
